# SeniorConnect_MasterLog.xlsx - "Auto update Excel log"
#
# Appends the latest batch of sensor-log rows (2026-02-01, the ~18:20-18:22
# bathroom-occupancy episode) to the bottom of five log sheets:
#   ALERTS, PIR, Humidity, Temperature, Proximity
#
# All log columns (Date/Timestamp/Hour/Value/...) are plain text in the
# source file, even when they look like a date ("2026-02-01"), a time
# ("18:21:21") or a percentage ("77.8%"). Excel's normal Range.Value setter
# auto-detects those patterns and silently stores them as a date serial /
# percentage number instead of text, so a helper is used that:
#   1. marks the cell's number format as Text ("@") before assigning,
#      which forces the literal string to be kept as typed;
#   2. clears that temporary format again right after, so the appended
#      cells end up with no explicit style - matching every other data
#      cell in these sheets (none of which carry an `s=` style).

function Set-SafeText {
    param($ws, $row, $col, $value)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- ALERTS: one new row (4) ---------------------------------------------
$ws = $wb.Worksheets.Item("ALERTS")
Set-SafeText $ws 4 1 '2026-02-01'
Set-SafeText $ws 4 2 '18:20:09'
Set-SafeText $ws 4 3 '18:00'
Set-SafeText $ws 4 4 'Bathroom'
Set-SafeText $ws 4 5 'MINIMAL'
Set-SafeText $ws 4 6 'MINIMAL ALERT: Bathroom occupied, no motion > 20s.'

# --- PIR: one new row (40) ------------------------------------------------
$ws = $wb.Worksheets.Item("PIR")
Set-SafeText $ws 40 1 '2026-02-01'
Set-SafeText $ws 40 2 '18:21:21'
Set-SafeText $ws 40 3 '18:00'
Set-SafeText $ws 40 4 'Bathroom'
Set-SafeText $ws 40 5 'Motion Detected'
Set-SafeText $ws 40 6 'Active'

# --- Humidity: eight new rows (37-44) -------------------------------------
$ws = $wb.Worksheets.Item("Humidity")
Set-SafeText $ws 37 1 '2026-02-01'
Set-SafeText $ws 37 2 '18:21:21'
Set-SafeText $ws 37 3 '18:00'
Set-SafeText $ws 37 4 'Bathroom'
Set-SafeText $ws 37 5 '77.8%'
Set-SafeText $ws 37 6 'Active'

Set-SafeText $ws 38 1 '2026-02-01'
Set-SafeText $ws 38 2 '18:21:26'
Set-SafeText $ws 38 3 '18:00'
Set-SafeText $ws 38 4 'Bathroom'
Set-SafeText $ws 38 5 '77.8%'
Set-SafeText $ws 38 6 'Active'

Set-SafeText $ws 39 1 '2026-02-01'
Set-SafeText $ws 39 2 '18:21:36'
Set-SafeText $ws 39 3 '18:00'
Set-SafeText $ws 39 4 'Bathroom'
Set-SafeText $ws 39 5 '78.0%'
Set-SafeText $ws 39 6 'Active'

Set-SafeText $ws 40 1 '2026-02-01'
Set-SafeText $ws 40 2 '18:21:41'
Set-SafeText $ws 40 3 '18:00'
Set-SafeText $ws 40 4 'Bathroom'
Set-SafeText $ws 40 5 '78.9%'
Set-SafeText $ws 40 6 'Active'

Set-SafeText $ws 41 1 '2026-02-01'
Set-SafeText $ws 41 2 '18:21:46'
Set-SafeText $ws 41 3 '18:00'
Set-SafeText $ws 41 4 'Bathroom'
Set-SafeText $ws 41 5 '77.7%'
Set-SafeText $ws 41 6 'Active'

Set-SafeText $ws 42 1 '2026-02-01'
Set-SafeText $ws 42 2 '18:21:51'
Set-SafeText $ws 42 3 '18:00'
Set-SafeText $ws 42 4 'Bathroom'
Set-SafeText $ws 42 5 '78.9%'
Set-SafeText $ws 42 6 'Active'

Set-SafeText $ws 43 1 '2026-02-01'
Set-SafeText $ws 43 2 '18:21:56'
Set-SafeText $ws 43 3 '18:00'
Set-SafeText $ws 43 4 'Bathroom'
Set-SafeText $ws 43 5 '77.9%'
Set-SafeText $ws 43 6 'Active'

Set-SafeText $ws 44 1 '2026-02-01'
Set-SafeText $ws 44 2 '18:22:01'
Set-SafeText $ws 44 3 '18:00'
Set-SafeText $ws 44 4 'Bathroom'
Set-SafeText $ws 44 5 '78.9%'
Set-SafeText $ws 44 6 'Active'

# --- Temperature: eight new rows (37-44) ----------------------------------
$ws = $wb.Worksheets.Item("Temperature")
Set-SafeText $ws 37 1 '2026-02-01'
Set-SafeText $ws 37 2 '18:21:22'
Set-SafeText $ws 37 3 '18:00'
Set-SafeText $ws 37 4 'Bathroom'
Set-SafeText $ws 37 5 '29.4C'
Set-SafeText $ws 37 6 'Active'

Set-SafeText $ws 38 1 '2026-02-01'
Set-SafeText $ws 38 2 '18:21:27'
Set-SafeText $ws 38 3 '18:00'
Set-SafeText $ws 38 4 'Bathroom'
Set-SafeText $ws 38 5 '29.4C'
Set-SafeText $ws 38 6 'Active'

Set-SafeText $ws 39 1 '2026-02-01'
Set-SafeText $ws 39 2 '18:21:37'
Set-SafeText $ws 39 3 '18:00'
Set-SafeText $ws 39 4 'Bathroom'
Set-SafeText $ws 39 5 '29.5C'
Set-SafeText $ws 39 6 'Active'

Set-SafeText $ws 40 1 '2026-02-01'
Set-SafeText $ws 40 2 '18:21:42'
Set-SafeText $ws 40 3 '18:00'
Set-SafeText $ws 40 4 'Bathroom'
Set-SafeText $ws 40 5 '29.4C'
Set-SafeText $ws 40 6 'Active'

Set-SafeText $ws 41 1 '2026-02-01'
Set-SafeText $ws 41 2 '18:21:47'
Set-SafeText $ws 41 3 '18:00'
Set-SafeText $ws 41 4 'Bathroom'
Set-SafeText $ws 41 5 '29.3C'
Set-SafeText $ws 41 6 'Active'

Set-SafeText $ws 42 1 '2026-02-01'
Set-SafeText $ws 42 2 '18:21:52'
Set-SafeText $ws 42 3 '18:00'
Set-SafeText $ws 42 4 'Bathroom'
Set-SafeText $ws 42 5 '29.4C'
Set-SafeText $ws 42 6 'Active'

Set-SafeText $ws 43 1 '2026-02-01'
Set-SafeText $ws 43 2 '18:21:57'
Set-SafeText $ws 43 3 '18:00'
Set-SafeText $ws 43 4 'Bathroom'
Set-SafeText $ws 43 5 '29.3C'
Set-SafeText $ws 43 6 'Active'

Set-SafeText $ws 44 1 '2026-02-01'
Set-SafeText $ws 44 2 '18:22:02'
Set-SafeText $ws 44 3 '18:00'
Set-SafeText $ws 44 4 'Bathroom'
Set-SafeText $ws 44 5 '29.4C'
Set-SafeText $ws 44 6 'Active'

# --- Proximity: four new rows (34-37) -------------------------------------
$ws = $wb.Worksheets.Item("Proximity")
Set-SafeText $ws 34 1 '2026-02-01'
Set-SafeText $ws 34 2 '18:20:09'
Set-SafeText $ws 34 3 '18:00'
Set-SafeText $ws 34 4 'Bathroom Door'
Set-SafeText $ws 34 5 'EXIT'
Set-SafeText $ws 34 6 'User EXITED Bathroom'

Set-SafeText $ws 35 1 '2026-02-01'
Set-SafeText $ws 35 2 '18:20:10'
Set-SafeText $ws 35 3 '18:00'
Set-SafeText $ws 35 4 'Bathroom Door'
Set-SafeText $ws 35 5 'ENTER'
Set-SafeText $ws 35 6 'User ENTERED Bathroom'

Set-SafeText $ws 36 1 '2026-02-01'
Set-SafeText $ws 36 2 '18:20:33'
Set-SafeText $ws 36 3 '18:00'
Set-SafeText $ws 36 4 'Bathroom Door'
Set-SafeText $ws 36 5 'EXIT'
Set-SafeText $ws 36 6 'User EXITED Bathroom'

Set-SafeText $ws 37 1 '2026-02-01'
Set-SafeText $ws 37 2 '18:20:47'
Set-SafeText $ws 37 3 '18:00'
Set-SafeText $ws 37 4 'Bathroom Door'
Set-SafeText $ws 37 5 'ENTER'
Set-SafeText $ws 37 6 'User ENTERED Bathroom'
